# dunn_matrix_equalw_multivar_norm_t189.xlsx
# "new sim results and new calculation"
#
# 1) Insert a new worksheet "sharpe_period" right before "VaR" - it starts
#    life as a full copy (values + formatting) of the existing "VaR" sheet,
#    then gets its own updated numbers (this is the old VaR matrix, recast
#    as the period (non-annualised) Sharpe ratio matrix).
# 2) Update the numeric matrices on the four pre-existing sheets
#    (annualised_return, mean_period_return, sharpe_annualized, VaR) with
#    the refreshed simulation results.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create "sharpe_period" as a duplicate of "VaR", inserted before it.
#    NOTE: after Worksheet.Copy() the original object reference can end
#    up bound to the newly-inserted copy (position-based COM binding), so
#    both sheets are re-fetched fresh by name afterwards.
# ---------------------------------------------------------------------
$origVar = $wb.Worksheets.Item("VaR")
$origVar.Copy($origVar)
$newCopy = $wb.Worksheets.Item("VaR (2)")
$newCopy.Name = "sharpe_period"

$sharpePeriod = $wb.Worksheets.Item("sharpe_period")
$varSheet = $wb.Worksheets.Item("VaR")

# ---------------------------------------------------------------------
# Helper data: [sheetName, cellRef, newValue] for every changed cell.
# ---------------------------------------------------------------------

# annualised_return
$annualisedReturn = $wb.Worksheets.Item("annualised_return")
$annualisedReturn.Range("G2").Value = 0.2241
$annualisedReturn.Range("D3").Value = 0
$annualisedReturn.Range("E3").Value = 0
$annualisedReturn.Range("H3").Value = 0
$annualisedReturn.Range("I3").Value = 0
$annualisedReturn.Range("J3").Value = 0
$annualisedReturn.Range("C4").Value = 0
$annualisedReturn.Range("H4").Value = 0
$annualisedReturn.Range("I4").Value = 0
$annualisedReturn.Range("J4").Value = 0
$annualisedReturn.Range("C5").Value = 0
$annualisedReturn.Range("H5").Value = 0
$annualisedReturn.Range("I5").Value = 0
$annualisedReturn.Range("J5").Value = 0
$annualisedReturn.Range("G6").Value = 0.9591
$annualisedReturn.Range("B7").Value = 0.2241
$annualisedReturn.Range("F7").Value = 0.9591
$annualisedReturn.Range("K7").Value = 0.0056
$annualisedReturn.Range("C8").Value = 0
$annualisedReturn.Range("D8").Value = 0
$annualisedReturn.Range("E8").Value = 0
$annualisedReturn.Range("C9").Value = 0
$annualisedReturn.Range("D9").Value = 0
$annualisedReturn.Range("E9").Value = 0
$annualisedReturn.Range("C10").Value = 0
$annualisedReturn.Range("D10").Value = 0
$annualisedReturn.Range("E10").Value = 0
$annualisedReturn.Range("G11").Value = 0.0056

# mean_period_return
$meanPeriodReturn = $wb.Worksheets.Item("mean_period_return")
$meanPeriodReturn.Range("G2").Value = 0.0198
$meanPeriodReturn.Range("D3").Value = 0
$meanPeriodReturn.Range("E3").Value = 0
$meanPeriodReturn.Range("H3").Value = 0.0026
$meanPeriodReturn.Range("I3").Value = 0
$meanPeriodReturn.Range("J3").Value = 0.0158
$meanPeriodReturn.Range("C4").Value = 0
$meanPeriodReturn.Range("H4").Value = 0
$meanPeriodReturn.Range("I4").Value = 0
$meanPeriodReturn.Range("J4").Value = 0
$meanPeriodReturn.Range("C5").Value = 0
$meanPeriodReturn.Range("H5").Value = 0
$meanPeriodReturn.Range("I5").Value = 0
$meanPeriodReturn.Range("J5").Value = 0
$meanPeriodReturn.Range("G6").Value = 0.8455
$meanPeriodReturn.Range("B7").Value = 0.0198
$meanPeriodReturn.Range("F7").Value = 0.8455
$meanPeriodReturn.Range("K7").Value = 0.0003
$meanPeriodReturn.Range("C8").Value = 0.0026
$meanPeriodReturn.Range("D8").Value = 0
$meanPeriodReturn.Range("E8").Value = 0
$meanPeriodReturn.Range("C9").Value = 0
$meanPeriodReturn.Range("D9").Value = 0
$meanPeriodReturn.Range("E9").Value = 0
$meanPeriodReturn.Range("C10").Value = 0.0158
$meanPeriodReturn.Range("D10").Value = 0
$meanPeriodReturn.Range("E10").Value = 0
$meanPeriodReturn.Range("G11").Value = 0.0003

# sharpe_annualized
$sharpeAnnualized = $wb.Worksheets.Item("sharpe_annualized")
$sharpeAnnualized.Range("D2").Value = 0.1055
$sharpeAnnualized.Range("E2").Value = 0.106
$sharpeAnnualized.Range("D3").Value = 0
$sharpeAnnualized.Range("E3").Value = 0
$sharpeAnnualized.Range("B4").Value = 0.1055
$sharpeAnnualized.Range("C4").Value = 0
$sharpeAnnualized.Range("B5").Value = 0.106
$sharpeAnnualized.Range("C5").Value = 0
$sharpeAnnualized.Range("G6").Value = 0.0017
$sharpeAnnualized.Range("F7").Value = 0.0017

# sharpe_period (new sheet - copy of old VaR, now updated)
$sharpePeriod.Range("D3").Value = 0.0003
$sharpePeriod.Range("E3").Value = 0.0003
$sharpePeriod.Range("C4").Value = 0.0003
$sharpePeriod.Range("C5").Value = 0.0003
$sharpePeriod.Range("G6").Value = 0.0006
$sharpePeriod.Range("F7").Value = 0.0006
$sharpePeriod.Range("J8").Value = 1
$sharpePeriod.Range("J9").Value = 1
$sharpePeriod.Range("H10").Value = 1
$sharpePeriod.Range("I10").Value = 1
$sharpePeriod.Range("K10").Value = 0
$sharpePeriod.Range("J11").Value = 0

# VaR (new matrix values for the existing sheet)
$varSheet.Range("D3").Value = 0
$varSheet.Range("E3").Value = 0
$varSheet.Range("F3").Value = 0.035
$varSheet.Range("G3").Value = 0.0008
$varSheet.Range("C4").Value = 0
$varSheet.Range("C5").Value = 0
$varSheet.Range("C6").Value = 0.035
$varSheet.Range("C7").Value = 0.0008
$varSheet.Range("J8").Value = 1
$varSheet.Range("J9").Value = 1
$varSheet.Range("H10").Value = 1
$varSheet.Range("I10").Value = 1
$varSheet.Range("K10").Value = 0
$varSheet.Range("J11").Value = 0
